$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Use an existing cell with style s="10" (e.g. I8) as a format template for
# the new "extra instructions" rows we are adding at rows 12 and 13.
$styleTemplate = $ws.Cells.Item(8, 9)

# --- Row 12: Load ---------------------------------------------------------
$styleTemplate.Copy($ws.Cells.Item(12, 9))
$styleTemplate.Copy($ws.Cells.Item(12, 10))
$styleTemplate.Copy($ws.Cells.Item(12, 11))
$styleTemplate.Copy($ws.Cells.Item(12, 12))
$styleTemplate.Copy($ws.Cells.Item(12, 13))

$ws.Cells.Item(12, 9).Value = "Load"
$ws.Cells.Item(12, 10).Value = "OpCode"
$ws.Cells.Item(12, 11).Value = "OpExt"
$ws.Cells.Item(12, 12).Value = "Address"
$ws.Cells.Item(12, 13).Value = "StoreTo"

# --- Row 13: Store ---------------------------------------------------------
$styleTemplate.Copy($ws.Cells.Item(13, 9))
$styleTemplate.Copy($ws.Cells.Item(13, 10))
$styleTemplate.Copy($ws.Cells.Item(13, 11))
$styleTemplate.Copy($ws.Cells.Item(13, 12))
$styleTemplate.Copy($ws.Cells.Item(13, 13))

$ws.Cells.Item(13, 9).Value = "Store"
$ws.Cells.Item(13, 10).Value = "OpCode"
$ws.Cells.Item(13, 11).Value = "OpExt"
$ws.Cells.Item(13, 12).Value = "Address"
$ws.Cells.Item(13, 13).Value = "StoreTo"

# Rows that now carry the extra I:M columns grow to the same row height used
# by the other annotated rows (8-11, 20, 22, 23) in this table.
$ws.Rows.Item(12).RowHeight = 14.9
$ws.Rows.Item(13).RowHeight = 14.9
